$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I4").Value = -0.7175322941186084
$ws.Range("J4").Value = 0.4690546999669646
$ws.Range("K4").Value = 0.6624895649182415
$ws.Range("L4").Value = 3.175250609064786
